$d = $word.ActiveDocument

$pairs = @(
    @("594÷9=", "786÷7="),
    @("151÷6=", "345÷2="),
    @("805÷6=", "466÷4="),
    @("541÷8=", "779÷7="),
    @("131÷5=", "999÷3="),
    @("293÷4=", "409÷6="),
    @("856÷4=", "909÷6="),
    @("160÷5=", "316÷9="),
    @("919÷8=", "546÷7="),
    @("535÷5=", "651÷2="),
    @("706÷5=", "197÷6="),
    @("109÷3=", "321÷3="),
    @("721÷2=", "689÷9="),
    @("152÷5=", "159÷8="),
    @("614÷5=", "346÷3="),
    @("740÷2=", "176÷2="),
    @("422÷2=", "177÷2="),
    @("445÷7=", "632÷6="),
    @("524÷6=", "156÷5="),
    @("142÷5=", "321÷9="),
    @("857÷2=", "962÷3="),
    @("876÷7=", "911÷8="),
    @("782÷7=", "319÷8="),
    @("943÷8=", "653÷7="),
    @("737÷8=", "659÷8=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
